# Remove the blank placeholder slide at position 12 (an empty slide with
# unused Text/Content/Title placeholders, sandwiched between the "Values of"
# slide and the "Example Problems!" slide). Deleting it shifts every
# subsequent slide up by one position; PowerPoint itself takes care of
# renumbering the relationship ids and refreshing any cached auto-update
# fields (date / slide-number placeholders) on save.
$p = $ppt.ActivePresentation
$p.Slides.Item(12).Delete()
